{"js": "const replacements = [\n  [\"88\u00d712=\", \"54\u00d730=\"],\n  [\"53\u00d728=\", \"55\u00d791=\"],\n  [\"24\u00d723=\", \"68\u00d798=\"],\n  [\"81\u00d774=\", \"52\u00d781=\"],\n  [\"11\u00d782=\", \"43\u00d756=\"],\n  [\"70\u00d741=\", \"70\u00d791=\"],\n  [\"39\u00d712=\", \"41\u00d761=\"],\n  [\"79\u00d728=\", \"91\u00d775=\"],\n  [\"21\u00d743=\", \"16\u00d750=\"],\n  [\"31\u00d717=\", \"48\u00d793=\"],\n  [\"72\u00d714=\", \"93\u00d782=\"],\n  [\"25\u00d717=\", \"89\u00d754=\"],\n  [\"16\u00d751=\", \"66\u00d774=\"],\n  [\"94\u00d718=\", \"45\u00d753=\"],\n  [\"45\u00d743=\", \"17\u00d756=\"],\n  [\"53\u00d723=\", \"64\u00d764=\"],\n  [\"20\u00d798=\", \"77\u00d771=\"],\n  [\"43\u00d796=\", \"87\u00d777=\"],\n  [\"87\u00d761=\", \"42\u00d715=\"],\n  [\"54\u00d771=\", \"82\u00d745=\"],\n  [\"35\u00d715=\", \"12\u00d743=\"],\n  [\"30\u00d795=\", \"18\u00d747=\"],\n  [\"79\u00d767=\", \"42\u00d776=\"],\n  [\"68\u00d795=\", \"44\u00d780=\"],\n  [\"51\u00d715=\", \"40\u00d716=\"],\n];\n\nconst body = context.document.body;\nfor (const [searchText, replaceText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${searchText}`);\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"88\u00d712=\", \"54\u00d730=\"),\n    @(\"53\u00d728=\", \"55\u00d791=\"),\n    @(\"24\u00d723=\", \"68\u00d798=\"),\n    @(\"81\u00d774=\", \"52\u00d781=\"),\n    @(\"11\u00d782=\", \"43\u00d756=\"),\n    @(\"70\u00d741=\", \"70\u00d791=\"),\n    @(\"39\u00d712=\", \"41\u00d761=\"),\n    @(\"79\u00d728=\", \"91\u00d775=\"),\n    @(\"21\u00d743=\", \"16\u00d750=\"),\n    @(\"31\u00d717=\", \"48\u00d793=\"),\n    @(\"72\u00d714=\", \"93\u00d782=\"),\n    @(\"25\u00d717=\", \"89\u00d754=\"),\n    @(\"16\u00d751=\", \"66\u00d774=\"),\n    @(\"94\u00d718=\", \"45\u00d753=\"),\n    @(\"45\u00d743=\", \"17\u00d756=\"),\n    @(\"53\u00d723=\", \"64\u00d764=\"),\n    @(\"20\u00d798=\", \"77\u00d771=\"),\n    @(\"43\u00d796=\", \"87\u00d777=\"),\n    @(\"87\u00d761=\", \"42\u00d715=\"),\n    @(\"54\u00d771=\", \"82\u00d745=\"),\n    @(\"35\u00d715=\", \"12\u00d743=\"),\n    @(\"30\u00d795=\", \"18\u00d747=\"),\n    @(\"79\u00d767=\", \"42\u00d776=\"),\n    @(\"68\u00d795=\", \"44\u00d780=\"),\n    @(\"51\u00d715=\", \"40\u00d716=\"),\n)\n\nforeach ($pair in $replacements) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $found = $find.Execute([ref]$searchText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$replaceText, 1)\n    if (-not $found) {\n        throw \"No match found for: $searchText\"\n    }\n}"}
